$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 344, pushing existing rows 344..451 down to 345..452.
# This duplicates the formatting/content of row 344 into the new row 344 first,
# then we overwrite the new row with its own data.
$ws.Rows.Item(344).Insert()

# Populate the freshly inserted row 344 with its data.
$ws.Cells.Item(344, 1).Value = 10
$ws.Cells.Item(344, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(344, 3).Value = "La Araucanía"
$ws.Cells.Item(344, 4).Value = 44988
$ws.Cells.Item(344, 5).Value = 9
$ws.Cells.Item(344, 6).Value = 100112044
$ws.Cells.Item(344, 7).Value = "Perejil"
$ws.Cells.Item(344, 8).Value = "Sin especificar"
$ws.Cells.Item(344, 9).Value = "Primera"
$ws.Cells.Item(344, 10).Value = 40
$ws.Cells.Item(344, 11).Value = 5000
$ws.Cells.Item(344, 12).Value = 5000
$ws.Cells.Item(344, 13).Value = 5000
$ws.Cells.Item(344, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(344, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(344, 16).Value = 1667
$ws.Cells.Item(344, 17).Value = 3
$ws.Cells.Item(344, 18).Value = "Hortaliza"
